{"js": "// Update the worksheet date header and the twenty-five \"two-digit \u00f7 one-digit\"\n// division prompts to the next day's generated values.\nconst replacements = [\n  [\"2026-01-26 Monday\", \"2026-01-27 Tuesday\"],\n  [\"40\u00f77=\", \"43\u00f79=\"],\n  [\"83\u00f76=\", \"50\u00f73=\"],\n  [\"87\u00f75=\", \"25\u00f73=\"],\n  [\"21\u00f72=\", \"97\u00f74=\"],\n  [\"10\u00f73=\", \"68\u00f75=\"],\n  [\"48\u00f78=\", \"80\u00f76=\"],\n  [\"55\u00f78=\", \"34\u00f79=\"],\n  [\"65\u00f77=\", \"33\u00f77=\"],\n  [\"79\u00f75=\", \"24\u00f72=\"],\n  [\"86\u00f75=\", \"94\u00f73=\"],\n  [\"65\u00f76=\", \"67\u00f73=\"],\n  [\"81\u00f72=\", \"32\u00f72=\"],\n  [\"19\u00f75=\", \"32\u00f75=\"],\n  [\"77\u00f76=\", \"12\u00f77=\"],\n  [\"22\u00f77=\", \"32\u00f75=\"],\n  [\"64\u00f76=\", \"76\u00f73=\"],\n  [\"63\u00f72=\", \"29\u00f79=\"],\n  [\"91\u00f75=\", \"20\u00f73=\"],\n  [\"18\u00f73=\", \"57\u00f75=\"],\n  [\"88\u00f72=\", \"27\u00f78=\"],\n  [\"33\u00f72=\", \"24\u00f73=\"],\n  [\"58\u00f78=\", \"88\u00f78=\"],\n  [\"78\u00f74=\", \"69\u00f75=\"],\n  [\"63\u00f78=\", \"24\u00f72=\"],\n  [\"27\u00f74=\", \"68\u00f74=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Could not find text to replace: ${oldText}`);\n  }\n\n  // Each source string is unique in the document, so only the first hit\n  // should exist; replace it (and, defensively, any further exact repeats).\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date header and the twenty-five \"two-digit \u00f7 one-digit\"\n# division prompts to the next day's generated values.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"2026-01-26 Monday\"; Replace = \"2026-01-27 Tuesday\" },\n    @{ Find = \"40\u00f77=\";  Replace = \"43\u00f79=\" },\n    @{ Find = \"83\u00f76=\";  Replace = \"50\u00f73=\" },\n    @{ Find = \"87\u00f75=\";  Replace = \"25\u00f73=\" },\n    @{ Find = \"21\u00f72=\";  Replace = \"97\u00f74=\" },\n    @{ Find = \"10\u00f73=\";  Replace = \"68\u00f75=\" },\n    @{ Find = \"48\u00f78=\";  Replace = \"80\u00f76=\" },\n    @{ Find = \"55\u00f78=\";  Replace = \"34\u00f79=\" },\n    @{ Find = \"65\u00f77=\";  Replace = \"33\u00f77=\" },\n    @{ Find = \"79\u00f75=\";  Replace = \"24\u00f72=\" },\n    @{ Find = \"86\u00f75=\";  Replace = \"94\u00f73=\" },\n    @{ Find = \"65\u00f76=\";  Replace = \"67\u00f73=\" },\n    @{ Find = \"81\u00f72=\";  Replace = \"32\u00f72=\" },\n    @{ Find = \"19\u00f75=\";  Replace = \"32\u00f75=\" },\n    @{ Find = \"77\u00f76=\";  Replace = \"12\u00f77=\" },\n    @{ Find = \"22\u00f77=\";  Replace = \"32\u00f75=\" },\n    @{ Find = \"64\u00f76=\";  Replace = \"76\u00f73=\" },\n    @{ Find = \"63\u00f72=\";  Replace = \"29\u00f79=\" },\n    @{ Find = \"91\u00f75=\";  Replace = \"20\u00f73=\" },\n    @{ Find = \"18\u00f73=\";  Replace = \"57\u00f75=\" },\n    @{ Find = \"88\u00f72=\";  Replace = \"27\u00f78=\" },\n    @{ Find = \"33\u00f72=\";  Replace = \"24\u00f73=\" },\n    @{ Find = \"58\u00f78=\";  Replace = \"88\u00f78=\" },\n    @{ Find = \"78\u00f74=\";  Replace = \"69\u00f75=\" },\n    @{ Find = \"63\u00f78=\";  Replace = \"24\u00f72=\" },\n    @{ Find = \"27\u00f74=\";  Replace = \"68\u00f74=\" }\n)\n\nforeach ($item in $replacements) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    # wdFindContinue = 1, wdReplaceAll = 2 \u2014 each Find string is unique in\n    # the document, so ReplaceAll is equivalent to a single substitution.\n    $rng.Find.Execute($item.Find, $false, $false, $false, $false, $false, $true, 1, $false, $item.Replace, 2) | Out-Null\n}\n"}
